$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 344.72726
$ws.Range("I15").Value = 344.72726
$ws.Range("K15").Value = 1034.18178
$ws.Range("M15").Value = -865.1817799999999
# Row 17
$ws.Range("H17").Value = 401114.7
$ws.Range("J17").Value = 401114.7
$ws.Range("L17").Value = 1203344.1
$ws.Range("N17").Value = -1203680.1
# Row 33
$ws.Range("H33").Value = 1149.091
$ws.Range("J33").Value = 2998.5
$ws.Range("L33").Value = 2998.5
$ws.Range("N33").Value = -3456.5
# Row 40
$ws.Range("H40").Value = 1665
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 1897.5
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 1897.5
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -2247.5
# Row 43
$ws.Range("H43").Value = 2399.5
$ws.Range("I43").Value = 3999
$ws.Range("J43").Value = 800
$ws.Range("K43").Value = 3999
$ws.Range("L43").Value = 800
$ws.Range("M43").Value = -3930
$ws.Range("N43").Value = -938
# Row 58
$ws.Range("H58").Value = 289.83334
$ws.Range("I58").Value = 289.83334
$ws.Range("K58").Value = 869.5000200000001
$ws.Range("M58").Value = -719.5000200000001
# Row 87
$ws.Range("H87").Value = 40329.168
$ws.Range("J87").Value = 40329.168
$ws.Range("L87").Value = 40329.168
$ws.Range("N87").Value = -42825.168
# Row 90
$ws.Range("H90").Value = 40329.168
$ws.Range("J90").Value = 40329.168
$ws.Range("L90").Value = 120987.504
$ws.Range("N90").Value = -133467.504
# Row 137
$ws.Range("H137").Value = 2125.0833
$ws.Range("I137").Value = 1565.4736
$ws.Range("K137").Value = 4696.4208
$ws.Range("M137").Value = -2146.4208

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 238.14285
$ws.Range("J5").Value = 299.33334
$ws.Range("L5").Value = 299.33334
$ws.Range("N5").Value = -523.33334
# Row 61
$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 8000
$ws.Range("K61").Value = 8000
$ws.Range("M61").Value = -7788
# Row 74
$ws.Range("H74").Value = 1013
$ws.Range("I74").Value = 974
$ws.Range("J74").Value = 1081.25
$ws.Range("K74").Value = 974
$ws.Range("L74").Value = 1081.25
$ws.Range("M74").Value = -100
$ws.Range("N74").Value = -2829.25
# Row 77
$ws.Range("H77").Value = 1013
$ws.Range("I77").Value = 974
$ws.Range("J77").Value = 1081.25
$ws.Range("K77").Value = 4870
$ws.Range("L77").Value = 5406.25
$ws.Range("M77").Value = -502
$ws.Range("N77").Value = -14142.25
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
# Row 136
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 8000
$ws.Range("K136").Value = 24000
$ws.Range("M136").Value = -21450
# Row 138
$ws.Range("H138").Value = 42060.5
$ws.Range("J138").Value = 42060.5
$ws.Range("L138").Value = 42060.5
$ws.Range("N138").Value = -52340.5
# Row 141
$ws.Range("H141").Value = 40429
$ws.Range("J141").Value = 40429
$ws.Range("L141").Value = 40429
$ws.Range("N141").Value = -50789

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 238.14285
$ws.Range("J4").Value = 299.33334
$ws.Range("L4").Value = 299.33334
$ws.Range("N4").Value = -529.33334

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 45.57143
$ws.Range("I7").Value = 44.833332
$ws.Range("K7").Value = 44.833332
$ws.Range("M7").Value = 68.166668

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2200
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2300
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 6900
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -8522
# Row 71
$ws.Range("H71").Value = 2200
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2300
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 20700
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -28812
# Row 128
$ws.Range("H128").Value = 324896.34
$ws.Range("I128").Value = 324896.34
$ws.Range("K128").Value = 974689.02
$ws.Range("M128").Value = -969709.02
# Row 131
$ws.Range("H131").Value = 670
$ws.Range("I131").Value = 670
$ws.Range("K131").Value = 2010
$ws.Range("M131").Value = 3030
# Row 139
$ws.Range("H139").Value = 2299.8
$ws.Range("I139").Value = 1624.75
$ws.Range("K139").Value = 4874.25
$ws.Range("M139").Value = 265.75

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2311
$ws.Range("J102").Value = 3998.5
$ws.Range("L102").Value = 3998.5
$ws.Range("N102").Value = -7242.5
# Row 132
$ws.Range("H132").Value = 1800
$ws.Range("I132").Value = 1800
$ws.Range("K132").Value = 5400
$ws.Range("M132").Value = -2870

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3566.5
$ws.Range("I46").Value = 3633.3333
$ws.Range("J46").Value = 3499.6667
$ws.Range("K46").Value = 3633.3333
$ws.Range("L46").Value = 3499.6667
$ws.Range("M46").Value = -3445.3333
$ws.Range("N46").Value = -3875.6667
# Row 47
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -15980
# Row 52
$ws.Range("H52").Value = 15000
$ws.Range("J52").Value = 15000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15466
# Row 61
$ws.Range("H61").Value = 3079.2727
$ws.Range("I61").Value = 3212.4285
$ws.Range("J61").Value = 2846.25
$ws.Range("K61").Value = 3212.4285
$ws.Range("L61").Value = 2846.25
$ws.Range("M61").Value = -3010.4285
$ws.Range("N61").Value = -3250.25
# Row 113
$ws.Range("H113").Value = 3079.2727
$ws.Range("I113").Value = 3212.4285
$ws.Range("J113").Value = 2846.25
$ws.Range("K113").Value = 3212.4285
$ws.Range("L113").Value = 2846.25
$ws.Range("M113").Value = -1042.4285
$ws.Range("N113").Value = -7186.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 11305
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = 22600
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = 22600
$ws.Range("M20").Value = 230
$ws.Range("N20").Value = -23080
# Row 62
$ws.Range("H62").Value = 6000
$ws.Range("J62").Value = 6000
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7248
# Row 65
$ws.Range("H65").Value = 6000
$ws.Range("J65").Value = 6000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240
# Row 136
$ws.Range("H136").Value = 3908.2942
$ws.Range("J136").Value = 3954.1428
$ws.Range("L136").Value = 11862.4284
$ws.Range("N136").Value = -16962.4284
